# Sari_Cases_3.xlsx template update
# - Renames a batch of header labels in row 1 (Spanish relabeling / typo fixes)
# - Gives the "Trabajo" header (S1) a top-vertical-aligned variant of its existing style
# - Selects A2 (just under the frozen header row)
# - Adjusts a few column widths to match the new header text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header relabels -------------------------------------------------
$ws.Range("S1").Value  = "Trabajo"
$ws.Range("W1").Value  = "TrabajoMercAnimVivos"
$ws.Range("X1").Value  = "País"
$ws.Range("Y1").Value  = "Departamento"
$ws.Range("Z1").Value  = "Municipio"
$ws.Range("AG1").Value = "IDPaísViaje.2Sem.Pre.Inici.Sint"
$ws.Range("AH1").Value = "IDArea.2Sem.Pre.Inici.Sint"
$ws.Range("AI1").Value = "IDEstado.2Sem.Pre.Inici.Sint"
$ws.Range("AJ1").Value = "IDBarrio.2Sem.Pre.Inici.Sint"
$ws.Range("AK1").Value = "Viaje.2SemPrevSintoma"
$ws.Range("BI1").Value = "Vac_Neumococo"
$ws.Range("DD1").Value = "Artralgia"
$ws.Range("DL1").Value = "Hepatomegalia"

# --- S1 gets a top-aligned variant of the existing bold/orange header style -
$ws.Range("S1").VerticalAlignment = -4160   # xlTop

# --- Selection moves to A2 (below the frozen header row) -------------------
$ws.Range("A2").Select()

# --- Column width touch-ups for the relabeled headers -----------------------
# (ColumnWidth is stored internally with a constant +5/6 character offset in
# this engine, so subtract it to land on the exact target character widths.)
$ws.Columns("Y").ColumnWidth  = 14 - 0.8333333333333334
$ws.Columns("AA").ColumnWidth = 16.5 - 0.8333333333333334
$ws.Columns("DK").ColumnWidth = 14.25 - 0.8333333333333334
